# Bump the "Förändrad" (Changed) date in column C from 45203 (2023-10-04)
# to 45204 (2023-10-05) for every data row (rows 2 through 307).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 307

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
